$d = $word.ActiveDocument

# Locate the standalone paragraph that contains only the italicized
# "2 Kings" text (the short-title line that immediately follows the
# "2KI" book-code heading) and remove the whole paragraph, including
# its paragraph mark, so the "2KI" heading paragraph is followed
# directly by the next (blank/space) paragraph.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if (($p.Range.Text.Trim() -eq "2 Kings") -and ($p.Range.Italic -eq -1)) {
        $p.Range.Delete()
        break
    }
}
